$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.914.95"
$ws.Range("E2").Value = "  -1.09%  "

$ws.Range("D3").Value = "1.925.42"
$ws.Range("E3").Value = "  +1.82%  "

$ws.Range("E4").Value = "  -0.16%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "320.18"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9993"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5062"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.26%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.4066"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.30%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.08357"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("E10").Value = "  -0.43%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.105"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "23.86"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.13%  "

$ws.Range("D13").Value = "1.921.43"
$ws.Range("E13").Value = "  +1.19%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.416"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.253"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "92.34"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.94%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001099"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.83%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06516"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.85%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "18.27"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.957"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").Value = "29.984.01"
$ws.Range("E23").Value = "  -0.83%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.33"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.55%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.189"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.75%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "22.16"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.74%  "

$ws.Range("D27").Value = "2.142.77"
$ws.Range("E27").Value = "  +1.44%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "162.17"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.21%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.343"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.73%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "129.22"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.133"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +4.35%  "

$ws.Range("E32").Value = "  -1.49%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.978"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.84%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.784"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.03%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.02451"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.41%  "

$ws.Range("E36").Value = "  +1.27%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06436"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.44%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2160"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.55%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.6531"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.77%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "8.757"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.196"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.76%  "

$ws.Range("E42").Value = "  -3.17%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.218"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.97%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.238"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +9.07%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.50"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.70%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.6097"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.35%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.608"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.212"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.93%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "122.20"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "79.11"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("E51").Value = "  -2.67%  "
